# Generate Report for Archive
#
# 1. Update status text "Ready for handoff" -> "In Translation" wherever it appears
#    (Overview!E2:F4, zh-cn!C2:C4, de-de!C2:C4).
# 2. Narrow the "Status" column(s) - Overview!E:F and zh-cn/de-de!C - from their
#    current ~17.22 character width down to ~13.41 character width.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        # NOTE: keep the string literal on the left of -eq. Some cells hold
        # the literal text "True"/"False" which Text/Value2 surface as a
        # native boolean; with the boolean on the left, PowerShell's -eq
        # would coerce the right-hand side to boolean too (any non-empty
        # string -> $true) and falsely match. Literal-left avoids that.
        if ("Ready for handoff" -eq $cell.Text) {
            $cell.Value = "In Translation"
        }
    }
}

# Note: ColumnWidth is quantized by the host to whole-pixel steps (~1/6 of a
# character unit here), so the nearest settable value that reproduces the
# target stored width (~13.4102) is obtained by assigning 12.5.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E1").EntireColumn.ColumnWidth = 12.5
$overview.Range("F1").EntireColumn.ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C1").EntireColumn.ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C1").EntireColumn.ColumnWidth = 12.5
